$wb = $excel.ActiveWorkbook

# --- Sheet: Pilot (sheet1) ---
$wsPilot = $wb.Worksheets.Item("Pilot")
$wsPilot.Activate()
$excel.ActiveWindow.ScrollRow = 20
$wsPilot.Range("K45").Select()

# --- Sheet: Leader (sheet2) ---
$wsLeader = $wb.Worksheets.Item("Leader")
$wsLeader.Activate()

# Update existing values (row 5)
$wsLeader.Range("N5").Value = 70
$wsLeader.Range("O5").Value = 75
$wsLeader.Range("Q5").Value = 80
$wsLeader.Range("S5").Value = 80

# Row 6
$wsLeader.Range("M6").Value = 75
$wsLeader.Range("N6").Value = 75
$wsLeader.Range("Q6").Value = 80
$wsLeader.Range("S6").Value = 80
$wsLeader.Range("Y6").Value = 80

# Row 7
$wsLeader.Range("N7").Value = 50
$wsLeader.Range("S7").Value = 60
$wsLeader.Range("V7").Value = 50
$wsLeader.Range("Y7").Value = 60

# Row 8
$wsLeader.Range("M8").Value = 70
$wsLeader.Range("N8").Value = 70
$wsLeader.Range("O8").Value = 70
$wsLeader.Range("Q8").Value = 75
$wsLeader.Range("R8").Value = 70
$wsLeader.Range("S8").Value = 75
$wsLeader.Range("Y8").Value = 75
$wsLeader.Range("AA8").Value = 70

# Row 9
$wsLeader.Range("G9").Value = 60
$wsLeader.Range("M9").Value = 60
$wsLeader.Range("N9").Value = 70
$wsLeader.Range("Q9").Value = 80
$wsLeader.Range("S9").Value = 80
$wsLeader.Range("Y9").Value = 80
$wsLeader.Range("AA9").Value = 60

# Row 10
$wsLeader.Range("F10").Value = 90
$wsLeader.Range("G10").Value = 85
$wsLeader.Range("M10").Value = 85
$wsLeader.Range("N10").Value = 75
$wsLeader.Range("O10").Value = 85
$wsLeader.Range("Q10").Value = 90
$wsLeader.Range("R10").Value = 85
$wsLeader.Range("Y10").Value = 90
$wsLeader.Range("AA10").Value = 85

# Row 11
$wsLeader.Range("G11").Value = 65
$wsLeader.Range("H11").Value = 65
$wsLeader.Range("J11").Value = 65
$wsLeader.Range("K11").Value = 65
$wsLeader.Range("L11").Value = 65
$wsLeader.Range("N11").Value = 65
$wsLeader.Range("P11").Value = 65
$wsLeader.Range("Q11").Value = 70
$wsLeader.Range("R11").Value = 65
$wsLeader.Range("V11").Value = 65
$wsLeader.Range("W11").Value = 65
$wsLeader.Range("Y11").Value = 70
$wsLeader.Range("Z11").Value = 65
$wsLeader.Range("AA11").Value = 65

# New row 14
$wsLeader.Range("F14").Value = 3
$wsLeader.Range("G14").Value = 2
$wsLeader.Range("H14").Value = 2
$wsLeader.Range("I14").Value = 2
$wsLeader.Range("J14").Value = 2
$wsLeader.Range("K14").Value = 2
$wsLeader.Range("L14").Value = 2
$wsLeader.Range("M14").Value = 2
$wsLeader.Range("N14").Value = 1
$wsLeader.Range("O14").Value = 2
$wsLeader.Range("P14").Value = 1
$wsLeader.Range("Q14").Value = 3
$wsLeader.Range("R14").Value = 2
$wsLeader.Range("S14").Value = 3
$wsLeader.Range("T14").Value = 1
$wsLeader.Range("U14").Value = 1
$wsLeader.Range("V14").Value = 1
$wsLeader.Range("W14").Value = 1
$wsLeader.Range("X14").Value = 1
$wsLeader.Range("Y14").Value = 3
$wsLeader.Range("Z14").Value = 1
$wsLeader.Range("AA14").Value = 2

# Apply style 2 (used by neighboring data cells) to the new row 14 cells
$wsLeader.Range("F14:AA14").Style = $wsLeader.Range("F12:AA12").Style

$wsLeader.Range("Q1").Select()
$excel.ActiveWindow.FreezePanes = $true
$wsLeader.Range("AA2").Select()

$wsPilot.Activate()
